$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.872.77'
$ws.Range("E2").Value = '  +2.41%  '

$ws.Range("D3").Value = '3.093.11'
$ws.Range("E3").Value = '  +5.48%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'579.85"
$ws.Range("E5").Value = '  +1.98%  '

$ws.Range("D6").Value = "'168.54"
$ws.Range("E6").Value = '  +6.39%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '3.088.53'
$ws.Range("E8").Value = '  +5.59%  '

$ws.Range("E9").Value = '  +1.34%  '

$ws.Range("E10").Value = '  -1.58%  '

$ws.Range("E11").Value = '  +4.26%  '

$ws.Range("E12").Value = '  +5.63%  '

$ws.Range("E13").Value = '  +2.20%  '

$ws.Range("D14").Value = "'36.42"
$ws.Range("E14").Value = '  +6.39%  '

$ws.Range("E15").Value = '  -0.53%  '

$ws.Range("D16").Value = '3.604.98'
$ws.Range("E16").Value = '  +5.39%  '

$ws.Range("D17").Value = '66.854.80'
$ws.Range("E17").Value = '  +2.32%  '

$ws.Range("E18").Value = '  +4.04%  '

$ws.Range("D19").Value = '3.092.94'
$ws.Range("E19").Value = '  +5.52%  '

$ws.Range("D20").Value = "'16.14"
$ws.Range("E20").Value = '  +4.87%  '

$ws.Range("D21").Value = "'467.44"
$ws.Range("E21").Value = '  +5.11%  '

$ws.Range("E22").Value = '  +4.12%  '

$ws.Range("E23").Value = '  +4.02%  '

$ws.Range("E24").Value = '  +2.06%  '

$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = '  +7.05%  '

$ws.Range("D26").Value = "'13.13"
$ws.Range("E26").Value = '  +8.61%  '

$ws.Range("D27").Value = "'10.16"
$ws.Range("E27").Value = '  +0.52%  '

$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("D30").Value = "'2.40"
$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("E31").Value = '  +4.13%  '

$ws.Range("D32").Value = "'0.0000103"
$ws.Range("E32").Value = '  +1.53%  '

$ws.Range("D33").Value = "'28.19"
$ws.Range("E33").Value = '  +4.19%  '

$ws.Range("E34").Value = '  +3.51%  '

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("E36").Value = '  +3.43%  '

$ws.Range("D37").Value = "'5.90"
$ws.Range("E37").Value = '  +3.14%  '

$ws.Range("D38").Value = "'2.11"
$ws.Range("E38").Value = '  +6.98%  '

$ws.Range("D39").Value = "'46.98"
$ws.Range("E39").Value = '  +4.70%  '

$ws.Range("D40").Value = "'50.33"
$ws.Range("E40").Value = '  +1.14%  '

$ws.Range("E41").Value = '  +6.53%  '

$ws.Range("E42").Value = '  +1.26%  '

$ws.Range("E43").Value = '  +2.63%  '

$ws.Range("E44").Value = '  -0.11%  '

$ws.Range("E45").Value = '  +3.07%  '

$ws.Range("D46").Value = "'383.07"
$ws.Range("E46").Value = '  -0.36%  '

$ws.Range("D47").Value = '2.779.76'

$ws.Range("D48").Value = "'135.18"
$ws.Range("E48").Value = '  +1.45%  '

$ws.Range("E49").Value = '  +0.05%  '

$ws.Range("D50").Value = "'25.02"
$ws.Range("E50").Value = '  +7.16%  '

$ws.Range("E51").Value = '  +2.33%  '
